$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the email addresses in column B (keep hyperlink style, just change the text)
$ws.Range("B2").Value = "user111223qqghr23344ffy6dedg@x.com"
$ws.Range("B3").Value = "user2345egdgfgfghrf4555ffff4e4j@l.com"
$ws.Range("B4").Value = "user32343948fdjj1234dfv5333o@g.com"
$ws.Range("B5").Value = "user498988hhjkhjbke33333eee3jj@a.com"

# Update the active cell selection to B3
$ws.Range("B3").Select()
